$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (prices/volume % stay as plain text, matching the
# original inlineStr cells). We force the Text number format before assigning so Excel
# does not auto-convert these numeric-looking / percent-looking strings into real numbers.
$updates = @{
    'D2' = '290.75'
    'E2' = '-4.18%'
    'D3' = '30.84'
    'E3' = '-5.91%'
    'D4' = '4.932'
    'E4' = '-0.20%'
    'D5' = '0.07201'
    'E5' = '-8.10%'
    'D6' = '1.796'
    'E6' = '-11.72%'
    'D7' = '7.668'
    'E7' = '-2.15%'
    'D8' = '3.747'
    'E8' = '-1.57%'
    'D9' = '0.8957'
    'E9' = '-3.07%'
    'D10' = '0.1655'
    'E10' = '-5.80%'
    'D11' = '0.07738'
    'E11' = '-1.22%'
    'D12' = '0.08056'
    'E12' = '-7.12%'
    'E13' = '-2.49%'
    'D14' = '0.1001'
    'E14' = '-0.33%'
    'D15' = '0.001494'
    'E15' = '-2.02%'
    'D16' = '0.005711'
    'E16' = '-3.83%'
    'D17' = '3.476'
    'E17' = '0.27%'
    'D18' = '2.080'
    'E18' = '-3.46%'
    'D20' = '0.1274'
    'E20' = '-1.26%'
    'D21' = '4.048'
    'E21' = '-6.26%'
    'D23' = '0.04513'
    'E23' = '-0.87%'
    'D24' = '0.001213'
    'E24' = '-0.94%'
    'D25' = '0.004007'
    'E25' = '-9.88%'
    'E26' = '-0.01%'
    'D39' = '0.01601'
    'E39' = '-7.98%'
    'D40' = '0.04386'
    'D41' = '0.007369'
    'E41' = '-1.42%'
    'D42' = '0.1309'
    'E42' = '-3.81%'
    'D43' = '0.007679'
    'D44' = '0.002071'
    'E44' = '-11.55%'
    'D45' = '0.009208'
    'E45' = '-12.96%'
    'D46' = '0.00005928'
    'E46' = '-4.81%'
    'E47' = '-0.05%'
    'D48' = '2.247'
    'E48' = '173.92%'
    'E50' = '-0.05%'
    'E51' = '-0.05%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
